$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.523.16"
$ws.Range("E2").Value = "  +0.52%  "

$ws.Range("D3").Value = "2.698.19"
$ws.Range("E3").Value = "  +2.11%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").Value = "'598.79"
$ws.Range("E5").Value = "  +0.25%  "

$ws.Range("D6").Value = "'160.27"
$ws.Range("E6").Value = "  +2.49%  "

$ws.Range("E7").Value = "  -0.09%  "

$ws.Range("E8").Value = "  +0.25%  "

$ws.Range("D9").Value = "2.696.32"
$ws.Range("E9").Value = "  +2.09%  "

$ws.Range("E10").Value = "  +0.36%  "

$ws.Range("E11").Value = "  -0.31%  "

$ws.Range("E12").Value = "  +0.92%  "

$ws.Range("E13").Value = "  +2.77%  "

$ws.Range("D14").Value = "'28.27"
$ws.Range("E14").Value = "  +1.11%  "

$ws.Range("D15").Value = "3.187.12"
$ws.Range("E15").Value = "  +1.98%  "

$ws.Range("E16").Value = "  -0.82%  "

$ws.Range("D17").Value = "68.439.23"
$ws.Range("E17").Value = "  +0.35%  "

$ws.Range("D18").Value = "2.715.43"
$ws.Range("E18").Value = "  +2.77%  "

$ws.Range("D19").Value = "'11.84"
$ws.Range("E19").Value = "  +4.17%  "

$ws.Range("D20").Value = "'364.97"
$ws.Range("E20").Value = "  +0.59%  "

$ws.Range("E21").Value = "  +4.07%  "

$ws.Range("E22").Value = "  +2.71%  "

$ws.Range("E23").Value = "  +2.20%  "

$ws.Range("E24").Value = "  +2.11%  "

$ws.Range("E25").Value = "  -1.47%  "

$ws.Range("E26").Value = "  -0.05%  "

$ws.Range("D27").Value = "'9.85"
$ws.Range("E27").Value = "  +1.11%  "

$ws.Range("E29").Value = "  +0.68%  "

$ws.Range("D30").Value = "'584.19"
$ws.Range("E30").Value = "  +5.20%  "

$ws.Range("D31").Value = "'0.998"
$ws.Range("E31").Value = "  -2.12%  "

$ws.Range("D32").Value = "'8.20"
$ws.Range("E32").Value = "  +1.97%  "

$ws.Range("E33").Value = "  +2.67%  "

$ws.Range("E34").Value = "  +5.15%  "

$ws.Range("E35").Value = "  +3.68%  "

$ws.Range("D36").Value = "'1.63"
$ws.Range("E36").Value = "  +6.06%  "

$ws.Range("E37").Value = "  -0.07%  "

$ws.Range("D38").Value = "'160.83"
$ws.Range("E38").Value = "  +0.14%  "

$ws.Range("E39").Value = "  +0.55%  "

$ws.Range("E40").Value = "  +2.16%  "

$ws.Range("E41").Value = "  +2.51%  "

$ws.Range("D42").Value = "'5.38"
$ws.Range("E42").Value = "  +1.08%  "

$ws.Range("B43").Value = "WhiteBITCoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D43").Value = "'17.98"
$ws.Range("E43").Value = "  +1.05%  "

$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").Value = "'2.67"
$ws.Range("E44").Value = "  +2.50%  "

$ws.Range("E45").Value = "  +0.04%  "

$ws.Range("E46").Value = "  -5.07%  "

$ws.Range("D47").Value = "'157.57"
$ws.Range("E47").Value = "  -0.64%  "

$ws.Range("E48").Value = "  +5.82%  "

$ws.Range("D49").Value = "'1.77"
$ws.Range("E49").Value = "  +5.07%  "

$ws.Range("D50").Value = "'0.602"
$ws.Range("E50").Value = "  +7.10%  "

$ws.Range("D51").Value = "'22.03"
$ws.Range("E51").Value = "  +0.47%  "
